# Update "想去人数" (number of people interested) values in column F
# for rows 2-5 on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1785
    $ws.Range("F3").Value = 8130
    $ws.Range("F4").Value = 188
    $ws.Range("F5").Value = 292
}
